$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.172.83"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.600.85"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "302.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3780"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.02"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.20%  "
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.266"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08115"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.591"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.412"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001244"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.49%  "
$ws.Range("D17").Value = "1.601.33"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06890"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.540"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").Value = "23.179.88"
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.400"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.977"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.255"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.371"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.747"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.90%  "
$ws.Range("D33").Value = "1.779.94"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9684"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07498"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2508"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08813"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.067"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7105"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.361"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6538"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.60%  "
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.015"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07952"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.201"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.208"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.45%  "
